$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 16
$ws.Range("H16").Value = 102814.66

# Row 18
$ws.Range("H18").Value = 7071.02

# Row 30
$ws.Range("F30").Value = 192.000
$ws.Range("H30").Value = 1191.06

# Row 32
$ws.Range("F32").Value = 1549.000
$ws.Range("H32").Value = 40347.18

# Row 35
$ws.Range("F35").Value = 878.000
$ws.Range("H35").Value = 120940.11

# Row 37
$ws.Range("F37").Value = 80.000
$ws.Range("H37").Value = 2968.87

# Row 41
$ws.Range("F41").Value = 3379.000
$ws.Range("H41").Value = 13063.11

# Row 49
$ws.Range("F49").Value = 49.000
$ws.Range("H49").Value = 1495.62

# Row 62
$ws.Range("H62").Value = 1776.62

# Row 64
$ws.Range("F64").Value = 183.000
$ws.Range("H64").Value = 35014.99

# Row 66
$ws.Range("H66").Value = 20870.80

# Row 86
$ws.Range("F86").Value = 939.850
$ws.Range("H86").Value = 9589.90

# Row 88
$ws.Range("F88").Value = 338.000
$ws.Range("H88").Value = 10618.76

$wb.Save()
